# Weekly price update: a new price-report row for "Albahaca" (Terminal La
# Palmera de La Serena) is inserted above the existing row 35, pushing the
# old rows 35-93 down to 36-94 (dimension grows from A1:R93 to A1:R94).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 35; Excel shifts rows 35..93 down to 36..94 and the
# new row inherits formatting (e.g. the date style in column D) from the
# row it displaces.
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row 35 with the new weekly observation.
$ws.Range("A35").Value = 8
$ws.Range("B35").Value = "Terminal La Palmera de La Serena"
$ws.Range("C35").Value = "Coquimbo"
$ws.Range("D35").Value = 44725
$ws.Range("E35").Value = 4
$ws.Range("F35").Value = 100112052
$ws.Range("G35").Value = "Albahaca"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 1000
$ws.Range("K35").Value = 3000
$ws.Range("L35").Value = 3500
$ws.Range("M35").Value = 3250
$ws.Range("N35").Value = "$/paquete"
$ws.Range("O35").Value = "Región de Arica y Parinacota"
$ws.Range("P35").Value = 3250
$ws.Range("Q35").Value = 1
$ws.Range("R35").Value = "Hortaliza"
